$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix delimiter typo in "Razon social" entries (comma -> period as separator)
$ws.Cells.Item(102, 5).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Cells.Item(176, 5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(188, 5).Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# Fix "Importe" column: remove thousands separators and normalize decimal
# separator to "." while keeping the values stored as text, as in the source data.
$ws.Range("H2:H220").NumberFormat = "@"

$ws.Cells.Item(2, 8).Value = "4200.00"
$ws.Cells.Item(3, 8).Value = "16611.00"
$ws.Cells.Item(4, 8).Value = "1605.00"
$ws.Cells.Item(5, 8).Value = "60000.00"
$ws.Cells.Item(6, 8).Value = "47827.50"
$ws.Cells.Item(7, 8).Value = "105000.00"
$ws.Cells.Item(8, 8).Value = "275800.00"
$ws.Cells.Item(9, 8).Value = "648000.00"
$ws.Cells.Item(10, 8).Value = "35000.00"
$ws.Cells.Item(11, 8).Value = "35000.00"
$ws.Cells.Item(12, 8).Value = "35000.00"
$ws.Cells.Item(13, 8).Value = "286800.00"
$ws.Cells.Item(14, 8).Value = "153400.00"
$ws.Cells.Item(15, 8).Value = "1322000.00"
$ws.Cells.Item(16, 8).Value = "6000.00"
$ws.Cells.Item(17, 8).Value = "29064.20"
$ws.Cells.Item(18, 8).Value = "4370.00"
$ws.Cells.Item(19, 8).Value = "286000.00"
$ws.Cells.Item(20, 8).Value = "462698.19"
$ws.Cells.Item(21, 8).Value = "135000.00"
$ws.Cells.Item(22, 8).Value = "3846.08"
$ws.Cells.Item(23, 8).Value = "2781.30"
$ws.Cells.Item(24, 8).Value = "1040.02"
$ws.Cells.Item(25, 8).Value = "400.00"
$ws.Cells.Item(26, 8).Value = "815786.19"
$ws.Cells.Item(27, 8).Value = "316980.30"
$ws.Cells.Item(28, 8).Value = "8256.00"
$ws.Cells.Item(29, 8).Value = "21.50"
$ws.Cells.Item(30, 8).Value = "209111.03"
$ws.Cells.Item(31, 8).Value = "6970.00"
$ws.Cells.Item(32, 8).Value = "48983.00"
$ws.Cells.Item(33, 8).Value = "73172.00"
$ws.Cells.Item(34, 8).Value = "56434.22"
$ws.Cells.Item(35, 8).Value = "26799.26"
$ws.Cells.Item(36, 8).Value = "21200.00"
$ws.Cells.Item(37, 8).Value = "3600.00"
$ws.Cells.Item(38, 8).Value = "26620.00"
$ws.Cells.Item(39, 8).Value = "4240.00"
$ws.Cells.Item(40, 8).Value = "11690.00"
$ws.Cells.Item(41, 8).Value = "12350.00"
$ws.Cells.Item(42, 8).Value = "837.30"
$ws.Cells.Item(43, 8).Value = "280.00"
$ws.Cells.Item(44, 8).Value = "1352.88"
$ws.Cells.Item(45, 8).Value = "39370.00"
$ws.Cells.Item(46, 8).Value = "27706.98"
$ws.Cells.Item(47, 8).Value = "220644.00"
$ws.Cells.Item(48, 8).Value = "894.75"
$ws.Cells.Item(49, 8).Value = "463.05"
$ws.Cells.Item(50, 8).Value = "66550.00"
$ws.Cells.Item(51, 8).Value = "432.00"
$ws.Cells.Item(52, 8).Value = "1073.44"
$ws.Cells.Item(53, 8).Value = "39165.77"
$ws.Cells.Item(54, 8).Value = "6938.40"
$ws.Cells.Item(55, 8).Value = "4597.00"
$ws.Cells.Item(56, 8).Value = "41816.35"
$ws.Cells.Item(57, 8).Value = "100.00"
$ws.Cells.Item(58, 8).Value = "10873.88"
$ws.Cells.Item(59, 8).Value = "550.00"
$ws.Cells.Item(60, 8).Value = "564.90"
$ws.Cells.Item(61, 8).Value = "44752.69"
$ws.Cells.Item(62, 8).Value = "26260.00"
$ws.Cells.Item(63, 8).Value = "4033.66"
$ws.Cells.Item(64, 8).Value = "682.44"
$ws.Cells.Item(65, 8).Value = "10.30"
$ws.Cells.Item(66, 8).Value = "277940.00"
$ws.Cells.Item(67, 8).Value = "244200.00"
$ws.Cells.Item(68, 8).Value = "12079.00"
$ws.Cells.Item(69, 8).Value = "9513.55"
$ws.Cells.Item(70, 8).Value = "1200.00"
$ws.Cells.Item(71, 8).Value = "3365.00"
$ws.Cells.Item(72, 8).Value = "72.13"
$ws.Cells.Item(73, 8).Value = "6836.63"
$ws.Cells.Item(74, 8).Value = "6922.50"
$ws.Cells.Item(75, 8).Value = "22119.98"
$ws.Cells.Item(76, 8).Value = "360.00"
$ws.Cells.Item(77, 8).Value = "17600.00"
$ws.Cells.Item(78, 8).Value = "20000.00"
$ws.Cells.Item(79, 8).Value = "815.00"
$ws.Cells.Item(80, 8).Value = "490.00"
$ws.Cells.Item(81, 8).Value = "8000.00"
$ws.Cells.Item(82, 8).Value = "41148.00"
$ws.Cells.Item(83, 8).Value = "5235.00"
$ws.Cells.Item(84, 8).Value = "3270.00"
$ws.Cells.Item(85, 8).Value = "22618.84"
$ws.Cells.Item(86, 8).Value = "450.00"
$ws.Cells.Item(87, 8).Value = "15400.00"
$ws.Cells.Item(88, 8).Value = "2700.00"
$ws.Cells.Item(89, 8).Value = "1584.72"
$ws.Cells.Item(90, 8).Value = "10000.00"
$ws.Cells.Item(91, 8).Value = "23552.00"
$ws.Cells.Item(92, 8).Value = "6052.00"
$ws.Cells.Item(93, 8).Value = "21000.00"
$ws.Cells.Item(94, 8).Value = "19920.00"
$ws.Cells.Item(95, 8).Value = "9355.00"
$ws.Cells.Item(96, 8).Value = "9960.00"
$ws.Cells.Item(97, 8).Value = "4900.00"
$ws.Cells.Item(98, 8).Value = "289.00"
$ws.Cells.Item(99, 8).Value = "196.38"
$ws.Cells.Item(100, 8).Value = "6520.00"
$ws.Cells.Item(101, 8).Value = "954.42"
$ws.Cells.Item(102, 8).Value = "500.00"
$ws.Cells.Item(103, 8).Value = "1095.00"
$ws.Cells.Item(104, 8).Value = "5998.00"
$ws.Cells.Item(105, 8).Value = "1908.00"
$ws.Cells.Item(106, 8).Value = "526097.59"
$ws.Cells.Item(107, 8).Value = "491.20"
$ws.Cells.Item(108, 8).Value = "56163.57"
$ws.Cells.Item(109, 8).Value = "88.67"
$ws.Cells.Item(110, 8).Value = "150.00"
$ws.Cells.Item(111, 8).Value = "75.00"
$ws.Cells.Item(112, 8).Value = "6354.00"
$ws.Cells.Item(113, 8).Value = "5484.00"
$ws.Cells.Item(114, 8).Value = "42000.00"
$ws.Cells.Item(115, 8).Value = "515.00"
$ws.Cells.Item(116, 8).Value = "488000.00"
$ws.Cells.Item(117, 8).Value = "5644.00"
$ws.Cells.Item(118, 8).Value = "381.00"
$ws.Cells.Item(119, 8).Value = "8000.00"
$ws.Cells.Item(120, 8).Value = "3178.50"
$ws.Cells.Item(121, 8).Value = "1540.00"
$ws.Cells.Item(122, 8).Value = "19620.00"
$ws.Cells.Item(123, 8).Value = "1837.50"
$ws.Cells.Item(124, 8).Value = "61649.85"
$ws.Cells.Item(125, 8).Value = "2297.00"
$ws.Cells.Item(126, 8).Value = "7958.26"
$ws.Cells.Item(127, 8).Value = "6354.98"
$ws.Cells.Item(128, 8).Value = "217.96"
$ws.Cells.Item(129, 8).Value = "1665.10"
$ws.Cells.Item(130, 8).Value = "614.61"
$ws.Cells.Item(131, 8).Value = "665.00"
$ws.Cells.Item(132, 8).Value = "18224.00"
$ws.Cells.Item(133, 8).Value = "6335.00"
$ws.Cells.Item(134, 8).Value = "3267.00"
$ws.Cells.Item(135, 8).Value = "4340.90"
$ws.Cells.Item(136, 8).Value = "7204.26"
$ws.Cells.Item(137, 8).Value = "7700.00"
$ws.Cells.Item(138, 8).Value = "45360.00"
$ws.Cells.Item(139, 8).Value = "6280.00"
$ws.Cells.Item(140, 8).Value = "410.00"
$ws.Cells.Item(141, 8).Value = "3500.00"
$ws.Cells.Item(142, 8).Value = "148800.00"
$ws.Cells.Item(143, 8).Value = "2236323.00"
$ws.Cells.Item(144, 8).Value = "50240.35"
$ws.Cells.Item(145, 8).Value = "6652.50"
$ws.Cells.Item(146, 8).Value = "5670.00"
$ws.Cells.Item(147, 8).Value = "1366222.00"
$ws.Cells.Item(148, 8).Value = "246895.00"
$ws.Cells.Item(149, 8).Value = "15000.00"
$ws.Cells.Item(150, 8).Value = "16000.00"
$ws.Cells.Item(151, 8).Value = "5000.00"
$ws.Cells.Item(152, 8).Value = "8000.00"
$ws.Cells.Item(153, 8).Value = "2500.00"
$ws.Cells.Item(154, 8).Value = "3867.50"
$ws.Cells.Item(155, 8).Value = "5112.00"
$ws.Cells.Item(156, 8).Value = "3000.00"
$ws.Cells.Item(157, 8).Value = "3000.00"
$ws.Cells.Item(158, 8).Value = "6000.00"
$ws.Cells.Item(159, 8).Value = "4000.00"
$ws.Cells.Item(160, 8).Value = "3000.00"
$ws.Cells.Item(161, 8).Value = "3000.00"
$ws.Cells.Item(162, 8).Value = "15000.00"
$ws.Cells.Item(163, 8).Value = "4000.00"
$ws.Cells.Item(164, 8).Value = "2500.00"
$ws.Cells.Item(165, 8).Value = "4100.00"
$ws.Cells.Item(166, 8).Value = "30000.00"
$ws.Cells.Item(167, 8).Value = "6000.00"
$ws.Cells.Item(168, 8).Value = "5000.00"
$ws.Cells.Item(169, 8).Value = "3000.00"
$ws.Cells.Item(170, 8).Value = "17880.00"
$ws.Cells.Item(171, 8).Value = "8000.00"
$ws.Cells.Item(172, 8).Value = "450.00"
$ws.Cells.Item(173, 8).Value = "2500.00"
$ws.Cells.Item(174, 8).Value = "1320.00"
$ws.Cells.Item(175, 8).Value = "224.19"
$ws.Cells.Item(176, 8).Value = "925.00"
$ws.Cells.Item(177, 8).Value = "258.16"
$ws.Cells.Item(178, 8).Value = "12150.00"
$ws.Cells.Item(179, 8).Value = "13245.00"
$ws.Cells.Item(180, 8).Value = "2778.03"
$ws.Cells.Item(181, 8).Value = "3482.66"
$ws.Cells.Item(182, 8).Value = "9782.00"
$ws.Cells.Item(183, 8).Value = "38900.00"
$ws.Cells.Item(184, 8).Value = "6016.96"
$ws.Cells.Item(185, 8).Value = "728.00"
$ws.Cells.Item(186, 8).Value = "6908.17"
$ws.Cells.Item(187, 8).Value = "3346.40"
$ws.Cells.Item(188, 8).Value = "2400.00"
$ws.Cells.Item(189, 8).Value = "1618.59"
$ws.Cells.Item(190, 8).Value = "10324.94"
$ws.Cells.Item(191, 8).Value = "5415.60"
$ws.Cells.Item(192, 8).Value = "502.16"
$ws.Cells.Item(193, 8).Value = "206.00"
$ws.Cells.Item(194, 8).Value = "1590.00"
$ws.Cells.Item(195, 8).Value = "2015.00"
$ws.Cells.Item(196, 8).Value = "4500.00"
$ws.Cells.Item(197, 8).Value = "16480.00"
$ws.Cells.Item(198, 8).Value = "6660.00"
$ws.Cells.Item(199, 8).Value = "4615.74"
$ws.Cells.Item(200, 8).Value = "330705.00"
$ws.Cells.Item(201, 8).Value = "18000.00"
$ws.Cells.Item(202, 8).Value = "2766.11"
$ws.Cells.Item(203, 8).Value = "5300.00"
$ws.Cells.Item(204, 8).Value = "6000.00"
$ws.Cells.Item(205, 8).Value = "21000.00"
$ws.Cells.Item(206, 8).Value = "6940092.41"
$ws.Cells.Item(207, 8).Value = "161000.00"
$ws.Cells.Item(208, 8).Value = "6000.00"
$ws.Cells.Item(209, 8).Value = "7754549.25"
$ws.Cells.Item(210, 8).Value = "215000.00"
$ws.Cells.Item(211, 8).Value = "6900.00"
$ws.Cells.Item(212, 8).Value = "5500.00"
$ws.Cells.Item(213, 8).Value = "4800.00"
$ws.Cells.Item(214, 8).Value = "6300.00"
$ws.Cells.Item(215, 8).Value = "23300.00"
$ws.Cells.Item(216, 8).Value = "8000.00"
$ws.Cells.Item(217, 8).Value = "127000.00"
$ws.Cells.Item(218, 8).Value = "4800.00"
$ws.Cells.Item(219, 8).Value = "68400.00"
$ws.Cells.Item(220, 8).Value = "5263.00"

# Restore the default (unformatted) style on the Importe column so only the
# underlying text changes, matching the original look of the sheet.
$ws.Range("H2:H220").Style = "Normal"
